$d = $word.ActiveDocument

# The last paragraph in the body (before the sectPr) currently ends with
# "Reading .json is easy, ... implement basic search/query first."
# Append three new paragraphs after it, each inheriting the same
# paragraph/run formatting (en-GB language) as the surrounding text.

$lastPara = $d.Paragraphs.Last
$r = $lastPara.Range

$p1 = "Ok very basic search/query from filename done! Used enter/submit instead of change detect due to performance issue from thumbnail reloading."
$p2 = "Now time to do some basic filter using metadata, so time to have a dictionary to metadata?"
$p3 = "Ok I think I need to do some more processing on  the extracted .json file for each week output as planned, then I can refer this to show filter drop down with number of item for each etc."

$r.InsertAfter("`r" + $p1 + "`r" + $p2 + "`r" + $p3)
